$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SS-read-write-data-json")

# Test case 1: Test Data value for row 11 changes from Test.json to Test0.json
$ws.Range("E11").Value = "Test0.json"

# Test case 2: "Name" field value gets appended with " / getStatisticsNames2"
$ws.Range("D19").Value = "getStatisticsNames1 / getStatisticsNames2"

# Test case 1: "Name" field value gets appended with " / getStatisticsFrom1"
$ws.Range("D5").Value = "getStatisticsFrom1 / getStatisticsFrom1"

# Update the active selection to D14 as in the diff
$ws.Activate()
$ws.Range("D14").Select()
